$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily refresh of the crypto price/volume table (GitHub Actions bot).
# For cells whose new text looks like a plain number (e.g. "112.37"), the
# cell is briefly stamped as Text ("@") before the assignment so Excel
# stores it as a string (matching the sheet's existing text-typed cells)
# instead of silently converting it to a numeric value; the format is
# then restored to General immediately after so no visible/style change
# is left behind.

$ws.Range("D2").Value = '43.203.02'
$ws.Range("E2").Value = '  -1.66%  '
$ws.Range("D3").Value = '2.275.63'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.37'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -3.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '265.02'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -2.03%  '
$ws.Range("E7").Value = '  -1.04%  '
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("E9").Value = '  -2.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.69'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("E11").Value = '  -1.82%  '
$ws.Range("E12").Value = '  -1.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.45'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -1.99%  '
$ws.Range("D15").Value = '2.616.53'
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.854'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '2.273.64'
$ws.Range("E17").Value = '  -1.67%  '
$ws.Range("D18").Value = '43.188.79'
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("E19").Value = '  -2.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.74'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.25'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.49'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -0.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.72'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("E24").Value = '  +1.91%  '
$ws.Range("E25").Value = '  -1.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.29'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("E28").Value = '  -1.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.30'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -5.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.32'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -3.47%  '
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '172.29'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -3.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.30'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0904'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -3.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.76'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +2.90%  '
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("E37").Value = '  -2.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0351'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -1.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.85'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -3.49%  '
$ws.Range("E40").Value = '  -6.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.62'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +9.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '76.58'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +7.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.75'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +7.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.236'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -3.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.09'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +0.82%  '
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("E47").Value = '  -2.34%  '
$ws.Range("E48").Value = '  -2.51%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.70'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.25'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0991'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -1.80%  '
